# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 5dcb06a3-... row on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-21 14:43:37"
$wsZh.Range("H3").Value = "2016-03-21 14:44:00"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-21 14:43:40"
$wsDe.Range("H3").Value = "2016-03-21 14:44:08"
